$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 4.077011333333332
$ws.Cells.Item(2, 8).Value = 12.231034
$ws.Cells.Item(2, 9).Value = 0.9715624748044627
$ws.Cells.Item(2, 10).Value = 0.9715624748044628
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.05601
$ws.Cells.Item(2, 14).Value = 0.16803
$ws.Cells.Item(2, 15).Value = 0.02710547761971223
$ws.Cells.Item(2, 16).Value = 0.02710547761971223
$ws.Cells.Item(2, 17).Value = 0.22835340478
$ws.Cells.Item(2, 18).Value = 2.05518064302
$ws.Cells.Item(2, 19).Value = 0.02633466491696459
$ws.Cells.Item(2, 20).Value = 0.0263346649169646

# Row 3 (ECs -> FAPs)
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 4.077011333333332
$ws.Cells.Item(3, 8).Value = 12.231034
$ws.Cells.Item(3, 9).Value = 0.9715624748044627
$ws.Cells.Item(3, 10).Value = 0.9715624748044628
$ws.Cells.Item(3, 14).Value = 5.594253
$ws.Cells.Item(3, 15).Value = 0.902427539668559
$ws.Cells.Item(3, 16).Value = 0.9024275396685592
$ws.Cells.Item(3, 17).Value = 7.602610960844665
$ws.Cells.Item(3, 18).Value = 68.42349864760199
$ws.Cells.Item(3, 19).Value = 0.8767647337720876
$ws.Cells.Item(3, 20).Value = 0.8767647337720879

# Row 4 (ECs -> MuSCs)
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 4.077011333333332
$ws.Cells.Item(4, 8).Value = 12.231034
$ws.Cells.Item(4, 9).Value = 0.9715624748044627
$ws.Cells.Item(4, 10).Value = 0.9715624748044628
$ws.Cells.Item(4, 13).Value = 0.145611
$ws.Cells.Item(4, 14).Value = 0.436833
$ws.Cells.Item(4, 15).Value = 0.07046698271172858
$ws.Cells.Item(4, 16).Value = 0.07046698271172858
$ws.Cells.Item(4, 17).Value = 0.5936576972579999
$ws.Cells.Item(4, 18).Value = 5.342919275321999
$ws.Cells.Item(4, 19).Value = 0.0684630761154103
$ws.Cells.Item(4, 20).Value = 0.0684630761154103

# Row 5 (FAPs -> ECs)
$ws.Cells.Item(5, 9).Value = 0.02843752519553723
$ws.Cells.Item(5, 10).Value = 0.02843752519553723
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.05601
$ws.Cells.Item(5, 14).Value = 0.16803
$ws.Cells.Item(5, 15).Value = 0.02710547761971223
$ws.Cells.Item(5, 16).Value = 0.02710547761971223
$ws.Cells.Item(5, 17).Value = 0.006683878670000001
$ws.Cells.Item(5, 18).Value = 0.06015490803000001
$ws.Cells.Item(5, 19).Value = 0.0007708127027476372
$ws.Cells.Item(5, 20).Value = 0.0007708127027476373

# Row 6 (FAPs -> FAPs)
$ws.Cells.Item(6, 9).Value = 0.02843752519553723
$ws.Cells.Item(6, 10).Value = 0.02843752519553723
$ws.Cells.Item(6, 14).Value = 5.594253
$ws.Cells.Item(6, 15).Value = 0.902427539668559
$ws.Cells.Item(6, 16).Value = 0.9024275396685592
$ws.Cells.Item(6, 19).Value = 0.02566280589647132
$ws.Cells.Item(6, 20).Value = 0.02566280589647133

# Row 7 (FAPs -> MuSCs)
$ws.Cells.Item(7, 9).Value = 0.02843752519553723
$ws.Cells.Item(7, 10).Value = 0.02843752519553723
$ws.Cells.Item(7, 13).Value = 0.145611
$ws.Cells.Item(7, 14).Value = 0.436833
$ws.Cells.Item(7, 15).Value = 0.07046698271172858
$ws.Cells.Item(7, 16).Value = 0.07046698271172858
$ws.Cells.Item(7, 17).Value = 0.017376294537
$ws.Cells.Item(7, 18).Value = 0.156386650833
$ws.Cells.Item(7, 19).Value = 0.002003906596318268
$ws.Cells.Item(7, 20).Value = 0.002003906596318268
